$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 210.92308
$ws.Range("I9").Value = 195.7
$ws.Range("K9").Value = 195.7
$ws.Range("M9").Value = -26.69999999999999
$ws.Range("H12").Value = 4441.2354
$ws.Range("I12").Value = 4325.1665
$ws.Range("J12").Value = 4504.5454
$ws.Range("K12").Value = 4325.1665
$ws.Range("L12").Value = 4504.5454
$ws.Range("M12").Value = -4155.1665
$ws.Range("N12").Value = -4844.5454
$ws.Range("H17").Value = 3383.9768
$ws.Range("J17").Value = 3557.775
$ws.Range("L17").Value = 10673.325
$ws.Range("N17").Value = -11009.325
$ws.Range("H33").Value = 837.25
$ws.Range("I33").Value = 837.25
$ws.Range("K33").Value = 837.25
$ws.Range("M33").Value = -608.25
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H112").Value = 5130.8965
$ws.Range("J112").Value = 5130.8965
$ws.Range("L112").Value = 15392.6895
$ws.Range("N112").Value = -17608.6895

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 664.0700000000001
$ws.Range("I32").Value = 590.5476
$ws.Range("K32").Value = 590.5476
$ws.Range("M32").Value = -303.5476
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51450
$ws.Range("H122").Value = 1099158
$ws.Range("I122").Value = 2800.1538
$ws.Range("J122").Value = 3474599.8
$ws.Range("K122").Value = 8400.4614
$ws.Range("L122").Value = 10423799.4
$ws.Range("M122").Value = -5950.4614
$ws.Range("N122").Value = -10428699.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 104997
$ws.Range("J59").Value = 104997
$ws.Range("L59").Value = 104997
$ws.Range("N59").Value = -106691
$ws.Range("H64").Value = 2249.25
$ws.Range("I64").Value = 2048.5
$ws.Range("J64").Value = 2450
$ws.Range("K64").Value = 2048.5
$ws.Range("L64").Value = 2450
$ws.Range("M64").Value = -1823.5
$ws.Range("N64").Value = -2900
$ws.Range("H67").Value = 2249.25
$ws.Range("I67").Value = 2048.5
$ws.Range("J67").Value = 2450
$ws.Range("K67").Value = 2048.5
$ws.Range("L67").Value = 2450
$ws.Range("M67").Value = -1268.5
$ws.Range("N67").Value = -4010
$ws.Range("H86").Value = 4353957.5
$ws.Range("I86").Value = 6674708.5
$ws.Range("J86").Value = 2549.625
$ws.Range("K86").Value = 6674708.5
$ws.Range("L86").Value = 2549.625
$ws.Range("M86").Value = -6673585.5
$ws.Range("N86").Value = -4795.625
$ws.Range("H89").Value = 4353957.5
$ws.Range("I89").Value = 6674708.5
$ws.Range("J89").Value = 2549.625
$ws.Range("K89").Value = 33373542.5
$ws.Range("L89").Value = 12748.125
$ws.Range("M89").Value = -33367926.5
$ws.Range("N89").Value = -23980.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1736.381
$ws.Range("I16").Value = 1384.2142
$ws.Range("J16").Value = 2440.7144
$ws.Range("K16").Value = 1384.2142
$ws.Range("L16").Value = 2440.7144
$ws.Range("M16").Value = -1097.2142
$ws.Range("N16").Value = -3014.7144
$ws.Range("H22").Value = 623.3333
$ws.Range("J22").Value = 870
$ws.Range("L22").Value = 870
$ws.Range("N22").Value = -1570
$ws.Range("H31").Value = 3704.32
$ws.Range("I31").Value = 1472.6666
$ws.Range("J31").Value = 4959.625
$ws.Range("K31").Value = 1472.6666
$ws.Range("L31").Value = 4959.625
$ws.Range("M31").Value = -1177.6666
$ws.Range("N31").Value = -5549.625
$ws.Range("H34").Value = 3704.32
$ws.Range("I34").Value = 1472.6666
$ws.Range("J34").Value = 4959.625
$ws.Range("K34").Value = 1472.6666
$ws.Range("L34").Value = 4959.625
$ws.Range("M34").Value = -1270.6666
$ws.Range("N34").Value = -5363.625
$ws.Range("H58").Value = 4169.3335
$ws.Range("I58").Value = 3835.6667
$ws.Range("J58").Value = 4503
$ws.Range("K58").Value = 3835.6667
$ws.Range("L58").Value = 4503
$ws.Range("M58").Value = -3632.6667
$ws.Range("N58").Value = -4909
$ws.Range("H94").Value = 1078.1875
$ws.Range("I94").Value = 763.6
$ws.Range("J94").Value = 1221.1818
$ws.Range("K94").Value = 763.6
$ws.Range("L94").Value = 1221.1818
$ws.Range("M94").Value = -312.6
$ws.Range("N94").Value = -2123.1818
$ws.Range("H113").Value = 1736.381
$ws.Range("I113").Value = 1384.2142
$ws.Range("J113").Value = 2440.7144
$ws.Range("K113").Value = 1384.2142
$ws.Range("L113").Value = 2440.7144
$ws.Range("M113").Value = 785.7858000000001
$ws.Range("N113").Value = -6780.7144
$ws.Range("H122").Value = 2951.1177
$ws.Range("I122").Value = 2977.6
$ws.Range("J122").Value = 2913.2856
$ws.Range("K122").Value = 8932.799999999999
$ws.Range("L122").Value = 8739.856800000001
$ws.Range("M122").Value = -6482.799999999999
$ws.Range("N122").Value = -13639.8568
$ws.Range("H134").Value = 25962.475
$ws.Range("I134").Value = 34006.395
$ws.Range("K134").Value = 102019.185
$ws.Range("M134").Value = -99484.185
$ws.Range("H136").Value = 4169.3335
$ws.Range("I136").Value = 3835.6667
$ws.Range("J136").Value = 4503
$ws.Range("K136").Value = 11507.0001
$ws.Range("L136").Value = 13509
$ws.Range("M136").Value = -8957.000100000001
$ws.Range("N136").Value = -18609
$ws.Range("H138").Value = 187999
$ws.Range("J138").Value = 187999
$ws.Range("L138").Value = 187999
$ws.Range("N138").Value = -198279

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 223385.86
$ws.Range("I46").Value = 334372.7
$ws.Range("J46").Value = 1412.2
$ws.Range("K46").Value = 1003118.1
$ws.Range("L46").Value = 4236.6
$ws.Range("M46").Value = -1003027.1
$ws.Range("N46").Value = -4418.6
$ws.Range("H56").Value = 10422783
$ws.Range("I56").Value = 10422783
$ws.Range("K56").Value = 10422783
$ws.Range("M56").Value = -10422253
$ws.Range("H113").Value = 2234.535
$ws.Range("I113").Value = 3709.1428
$ws.Range("J113").Value = 1522.6552
$ws.Range("K113").Value = 11127.4284
$ws.Range("L113").Value = 4567.9656
$ws.Range("M113").Value = -8957.428400000001
$ws.Range("N113").Value = -8907.9656
$ws.Range("H122").Value = 1074.4546
$ws.Range("I122").Value = 1018
$ws.Range("J122").Value = 1121.5
$ws.Range("K122").Value = 9162
$ws.Range("L122").Value = 10093.5
$ws.Range("M122").Value = -6712
$ws.Range("N122").Value = -14993.5
$ws.Range("H139").Value = 41668824
$ws.Range("I139").Value = 62501820
$ws.Range("K139").Value = 187505460
$ws.Range("M139").Value = -187500320

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 55557224
$ws.Range("I113").Value = 166666670
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 166666670
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = -166664500
$ws.Range("N113").Value = -6839.5
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178
$ws.Range("H132").Value = 3821.7896
$ws.Range("I132").Value = 3170.3845
$ws.Range("J132").Value = 5233.1665
$ws.Range("K132").Value = 9511.1535
$ws.Range("L132").Value = 15699.4995
$ws.Range("M132").Value = -6981.1535
$ws.Range("N132").Value = -20759.4995
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4646.1665
$ws.Range("I40").Value = 3795.318
$ws.Range("K40").Value = 3795.318
$ws.Range("M40").Value = -3659.318
$ws.Range("H46").Value = 4999.2188
$ws.Range("J46").Value = 7645.273
$ws.Range("L46").Value = 7645.273
$ws.Range("N46").Value = -8021.273
$ws.Range("H61").Value = 5851658.5
$ws.Range("I61").Value = 6948301
$ws.Range("K61").Value = 6948301
$ws.Range("M61").Value = -6948099
$ws.Range("H113").Value = 5851658.5
$ws.Range("I113").Value = 6948301
$ws.Range("K113").Value = 6948301
$ws.Range("M113").Value = -6946131
$ws.Range("H132").Value = 4836.8677
$ws.Range("I132").Value = 4703.5103
$ws.Range("K132").Value = 14110.5309
$ws.Range("M132").Value = -11580.5309

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 16006
$ws.Range("J31").Value = 21509
$ws.Range("L31").Value = 21509
$ws.Range("N31").Value = -22205
$ws.Range("H136").Value = 2002.2878
$ws.Range("I136").Value = 978.4091
$ws.Range("K136").Value = 2935.2273
$ws.Range("M136").Value = -385.2273
